# The Pearson logo picture (embedded in the default + first-page footers)
# was renamed from image1.png to image2.png, and the BTEC logo picture
# (embedded in the default + first-page headers) was renamed from
# image2.jpg to image1.jpg.
#
# InlineShape has no settable "Name" in the Word object model, so each
# picture is briefly promoted to a floating Shape (which does expose
# .Name), renamed, then converted back to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Footers: default (Item 1) and first-page (Item 2) both carry the
# Pearson logo - rename image1.png -> image2.png in both.
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Headers: default (Item 1) and first-page (Item 2) both carry the
# BTEC logo - rename image2.jpg -> image1.jpg in both.
Rename-InlinePicture $sec.Headers.Item(1).Range "image1.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
